$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Refresh the "datetimeFigureOut" date placeholders (slide master + every
#    slide layout) from 5/13/2022 -> 5/15/2022, as happens automatically when
#    the deck is reopened/saved on a later day.
# ---------------------------------------------------------------------------
$oldDate = "5/13/2022"
$newDate = "5/15/2022"

$targets = @($p.SlideMaster)
for ($i = 1; $i -le $p.SlideMaster.CustomLayouts.Count; $i++) {
    $targets += $p.SlideMaster.CustomLayouts.Item($i)
}

foreach ($t in $targets) {
    for ($j = 1; $j -le $t.Shapes.Count; $j++) {
        $sh = $t.Shapes.Item($j)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# ---------------------------------------------------------------------------
# 2) DFA drawing fix on slide 2: the transition label between the two middle
#    states was mislabeled "i" - it should read "l". Target the shape by its
#    stable shape Id (Group 32 -> TextBox 34) so we only ever touch that one
#    run, regardless of any other single-letter labels on the drawing.
# ---------------------------------------------------------------------------
$s = $p.Slides.Item(2)
for ($k = 1; $k -le $s.Shapes.Count; $k++) {
    $top = $s.Shapes.Item($k)
    if ($top.Id -eq 33) {
        for ($m = 1; $m -le $top.GroupItems.Count; $m++) {
            $inner = $top.GroupItems.Item($m)
            if ($inner.Id -eq 35) {
                $inner.TextFrame.TextRange.Text = "l"
            }
        }
    }
}
